$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Low frequency block (rows 17-32) ---
$ws.Range("D18").Value = "Mean increase"
$ws.Range("F18").Value = "Median increase"
$ws.Range("D18").Font.Bold = $true
$ws.Range("F18").Font.Bold = $true

$ws.Range("D19").Formula = "= ((E3 / 95.321842) * 100) - 100"
$ws.Range("F19").Formula = "= ((E10 / 95.22216) * 100) - 100"
$ws.Range("D19").ClearFormats()
$ws.Range("F19").ClearFormats()

# --- Medium frequency block (rows 49-64) ---
$ws.Range("D50").Value = "Mean increase"
$ws.Range("F50").Value = "Median increase"
$ws.Range("D50").Font.Bold = $true
$ws.Range("F50").Font.Bold = $true

$ws.Range("D51").Formula = "= ((E35 / 95.321842) * 100) - 100"
$ws.Range("F51").Formula = "= ((E42 / 95.22216) * 100) - 100"
$ws.Range("D51").ClearFormats()
$ws.Range("F51").ClearFormats()

# --- High frequency block (rows 81-96) ---
$ws.Range("D82").Value = "Mean increase"
$ws.Range("F82").Value = "Median increase"
$ws.Range("D82").Font.Bold = $true
$ws.Range("F82").Font.Bold = $true

$ws.Range("D83").Formula = "=((E67 / 95.321842) * 100) - 100"
$ws.Range("F83").Formula = "= ((E74 / 95.22216) * 100) - 100"
$ws.Range("D83").ClearFormats()
$ws.Range("F83").ClearFormats()

# --- All (combined) block (rows 113-114) ---
$ws.Range("D113").Value = "Mean increase"
$ws.Range("F113").Value = "Median increase"
$ws.Range("D113").Font.Bold = $true
$ws.Range("F113").Font.Bold = $true

$ws.Range("D114").Formula = "= (D19 + D51 + D83) / 3"
$ws.Range("F114").Formula = "= (F19 + F51 + F83) / 3"
$ws.Range("D114").ClearFormats()
$ws.Range("F114").ClearFormats()

# --- View / window adjustments ---
$ws.Application.ActiveWindow.ScrollRow = 87
$ws.Range("H113").Select()
